$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "30.227.94"
$ws.Range("E2").Value = "  +6.36%  "
$ws.Range("D3").Value = "1.920.82"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").Value = "'330.78"
$ws.Range("E5").Value = "  +5.26%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("D8").Value = "'0.4091"
$ws.Range("E8").Value = "  +5.23%  "
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "'42.99"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "'1.129"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "'23.12"
$ws.Range("E12").Value = "  +14.16%  "
$ws.Range("D13").Value = "'6.469"
$ws.Range("E13").Value = "  +4.88%  "
$ws.Range("D14").Value = "1.903.40"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "'7.416"
$ws.Range("E15").Value = "  +2.69%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "'95.33"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "'0.06705"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'18.56"
$ws.Range("E20").Value = "  +5.80%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'6.034"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").Value = "30.230.54"
$ws.Range("E23").Value = "  +6.21%  "
$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").Value = "'2.220"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "2.126.22"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "'161.59"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").Value = "'21.18"
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("D29").Value = "'2.424"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "'129.24"
$ws.Range("E30").Value = "  +3.00%  "
$ws.Range("D31").Value = "'1.096"
$ws.Range("E31").Value = "  +6.22%  "
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("D33").Value = "'6.028"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'3.607"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").Value = "'0.02499"
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "'0.06596"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "'0.2220"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").Value = "'1.233"
$ws.Range("E38").Value = "  +4.81%  "
$ws.Range("D39").Value = "'5.187"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.88"
$ws.Range("E40").Value = "  +7.41%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.833"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").Value = "'0.6546"
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").Value = "'1.244"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "'0.6184"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").Value = "'3.751"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").Value = "'2.104"
$ws.Range("E47").Value = "  +5.83%  "
$ws.Range("D48").Value = "'1.247"
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("D49").Value = "'124.32"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "'1.167"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'79.78"
$ws.Range("E51").Value = "  +5.42%  "
